$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 12:22"

# Reorder Rumania ahead of Catar (row 35 becomes Rumania, row 36 becomes Catar)
# and refresh both countries' COVID figures.
$ws.Range("A35").Value = "Rumania"
$ws.Range("B35").Value = 11978
$ws.Range("C35").Value = 362
$ws.Range("D35").Value = 3569
$ws.Range("E35").Value = 7734
$ws.Range("F35").Value = 247
$ws.Range("G35").Value = 12
$ws.Range("H35").Value = 675

$ws.Range("A36").Value = "Catar"
$ws.Range("B36").Value = 11921
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 1134
$ws.Range("E36").Value = 10777
$ws.Range("F36").Value = 72
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 10

# Marruecos (row 55) - updated figures
$ws.Range("B55").Value = 4289
$ws.Range("C55").Value = 37
$ws.Range("D55").Value = 890
$ws.Range("E55").Value = 3232
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 167

# Uzbekistan (row 68) - updated figures
$ws.Range("D68").Value = 1055
$ws.Range("E68").Value = 892

# Hong Kong (row 90) - updated figures
$ws.Range("D90").Value = 830
$ws.Range("E90").Value = 204
